# Daily attendance processing - 2025-11-14 18:29:47
# Normalize the "Recorded By" (column G) values: move the trailing
# contributor to the front of the comma-separated list, for every row
# whose current value matches one of the known legacy orderings.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$map = @{
    "System, dnasr281@gmail.com"          = "dnasr281@gmail.com, System"
    "admin@admin.com, dnasr281@gmail.com" = "dnasr281@gmail.com, admin@admin.com"
    "System, backup@backdoor.com, system" = "system, System, backup@backdoor.com"
}

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1

for ($row = 1; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)   # Column G = "Recorded By"
    $value = $cell.Value2

    if ($null -ne $value -and $map.ContainsKey($value)) {
        $cell.Value2 = $map[$value]
    }
}
